$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the URL path text for the search endpoints (remove the {title}/{label} path params)
$ws.Range("A2").Value = "notes/search_title"
$ws.Range("A3").Value = "notes/search_label"

# Update the HTTP method for those two endpoints from GET to POST
$ws.Range("B2").Value = "POST"
$ws.Range("B3").Value = "POST"

# Update the selected cell in the sheet view
$ws.Range("D4").Select()
